# Prepare to integrate SDK
# Insert a new error-code row ("APPUSER_CHECK_NOT_PASS" / 12001) into the
# "错误码" (error code) sheet, just above the RESPONSE_TIME_OUT row, and
# select/scroll the sheet roughly where the author left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A blank row (43) is kept before the new entry and a blank row (45) is kept
# after it, matching the sheet's existing layout convention of separating
# numbered "blocks" of error codes with an empty row. Inserting a single row
# at row 44 pushes the existing RESPONSE_TIME_OUT row (was row 45) down to
# row 46, opening up row 44 for the new entry.
$ws.Rows.Item(44).Insert()

$ws.Range("A44").Value = "APPUSER_CHECK_NOT_PASS"
$ws.Range("B44").Value = 12001

# Leave the selection/scroll roughly where the editor was working (near the
# bottom of the table where the new row was added).
$ws.Range("A1:XFD1048576").Select()
$ws.Range("F24").Activate()
